$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update election results row (row 2) with new vote counts per party/column.
$ws.Range("H2").Value = 22
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 320
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 85
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 55
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 27
$ws.Range("T2").Value = 50
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 482
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 480
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 3
$ws.Range("AA2").Value = 0
